$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.22914666016441
$ws.Range("C2").Value = 0.22914666016441
$ws.Range("D2").Value = 1.45647640011177
$ws.Range("E2").Value = 0.00784726815629703
$ws.Range("F2").Value = 0.183

# Row 3
$ws.Range("B3").Value = 0.663177523997795
$ws.Range("C3").Value = 0.663177523997795
$ws.Range("D3").Value = 4.21521488506236
$ws.Range("E3").Value = 0.0227109217402772
$ws.Range("F3").Value = 0.0006

# Row 4
$ws.Range("B4").Value = 5.7963980110554
$ws.Range("C4").Value = 1.9321326703518
$ws.Range("D4").Value = 12.2808058133304
$ws.Range("E4").Value = 0.198501210974416

# Row 5
$ws.Range("B5").Value = 0.120936833529588
$ws.Range("C5").Value = 0.120936833529588
$ws.Range("D5").Value = 0.768685189710867
$ws.Range("E5").Value = 0.00414155616319791
$ws.Range("F5").Value = 0.6117

# Row 6
$ws.Range("B6").Value = 0.297482721969514
$ws.Range("C6").Value = 0.0991609073231713
$ws.Range("D6").Value = 0.630275480455382
$ws.Range("E6").Value = 0.0101874785758824
$ws.Range("F6").Value = 0.8889

# Row 7
$ws.Range("B7").Value = 1.10879878102872
$ws.Range("C7").Value = 0.369599593676241
$ws.Range("D7").Value = 2.34920764410929
$ws.Range("E7").Value = 0.0379714954600026
$ws.Range("F7").Value = 0.0007

# Row 8
$ws.Range("B8").Value = 0.217388815655829
$ws.Range("C8").Value = 0.0724629385519429
$ws.Range("D8").Value = 0.460580834160665
$ws.Range("E8").Value = 0.00744461354752951
$ws.Range("F8").Value = 0.9766

# Row 9
$ws.Range("B9").Value = 20.7674900461009
$ws.Range("C9").Value = 0.157329470046219
$ws.Range("E9").Value = 0.711195455382397

# Row 10
$ws.Range("B10").Value = 29.2008193935022
